$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$sh = $m.Shapes.Item(18)

Write-Host "--- Left ---"
try { $sh.Left = 10; Write-Host "OK" } catch { Write-Host "FAIL: $($_.Exception.Message)" }
Write-Host "--- Top ---"
try { $sh.Top = 10; Write-Host "OK" } catch { Write-Host "FAIL: $($_.Exception.Message)" }
Write-Host "--- Width ---"
try { $sh.Width = 10; Write-Host "OK" } catch { Write-Host "FAIL: $($_.Exception.Message)" }
Write-Host "--- Height ---"
try { $sh.Height = 10; Write-Host "OK" } catch { Write-Host "FAIL: $($_.Exception.Message)" }
Write-Host "--- Rotation ---"
try { $sh.Rotation = 10; Write-Host "OK" } catch { Write-Host "FAIL: $($_.Exception.Message)" }
Write-Host "--- Name ---"
try { $sh.Name = "x"; Write-Host "OK" } catch { Write-Host "FAIL: $($_.Exception.Message)" }
Write-Host "--- Fill.ForeColor ---"
try { $sh.Fill.ForeColor.RGB = 255; Write-Host "OK" } catch { Write-Host "FAIL: $($_.Exception.Message)" }
Write-Host "--- Line.Visible ---"
try { $sh.Line.Visible = $true; Write-Host "OK" } catch { Write-Host "FAIL: $($_.Exception.Message)" }
Write-Host "--- TextFrame.TextRange.Text ---"
try { $sh.TextFrame.TextRange.Text = "abc"; Write-Host "OK" } catch { Write-Host "FAIL: $($_.Exception.Message)" }
Write-Host "--- TextFrame.TextRange.Font.Bold ---"
try { $sh.TextFrame.TextRange.Font.Bold = $true; Write-Host "OK" } catch { Write-Host "FAIL: $($_.Exception.Message)" }
Write-Host "--- Visible ---"
try { $sh.Visible = $true; Write-Host "OK" } catch { Write-Host "FAIL: $($_.Exception.Message)" }
